$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26, shifting existing rows 26..145 down to 27..146
$ws.Rows.Item(26).Insert()

# Populate the new row 26 with the new data
$ws.Cells.Item(26, 1).Value = 10
$ws.Cells.Item(26, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(26, 3).Value = "La Araucanía"
$ws.Cells.Item(26, 4).Value = 44707
$ws.Cells.Item(26, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(26, 5).Value = 9
$ws.Cells.Item(26, 6).Value = 100112012
$ws.Cells.Item(26, 7).Value = "Espinaca"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 50
$ws.Cells.Item(26, 11).Value = 9000
$ws.Cells.Item(26, 12).Value = 9000
$ws.Cells.Item(26, 13).Value = 9000
$ws.Cells.Item(26, 14).Value = "`$/docena de atados"
$ws.Cells.Item(26, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(26, 16).Value = 3000
$ws.Cells.Item(26, 17).Value = 3
$ws.Cells.Item(26, 18).Value = "Hortaliza"
